$d = $word.ActiveDocument
$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

# --- New list paragraph: "Change word default format to doc 97/2003/xp" ---
$last = $d.Paragraphs.Last.Range
$last.InsertParagraphAfter()

$p2 = $d.Paragraphs.Last.Range
$p2xml = '<w:p xmlns:w="' + $w + '">' +
  '<w:pPr>' +
    '<w:pStyle w:val="ListParagraph"/>' +
    '<w:numPr><w:ilvl w:val="0"/><w:numId w:val="1"/></w:numPr>' +
    '<w:rPr><w:lang w:val="en-US"/></w:rPr>' +
  '</w:pPr>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>Change word default format to doc 97/2003/</w:t></w:r>' +
  '<w:proofErr w:type="spellStart"/>' +
  '<w:r><w:rPr><w:lang w:val="en-US"/></w:rPr><w:t>xp</w:t></w:r>' +
  '<w:proofErr w:type="spellEnd"/>' +
  '</w:p>'
$p2.InsertXML($p2xml)

# --- Two trailing empty paragraphs (plain, no list formatting) ---
$p2again = $d.Paragraphs.Last.Range
$p2again.InsertParagraphAfter()
$p3 = $d.Paragraphs.Last.Range
$p3xml = '<w:p xmlns:w="' + $w + '"><w:pPr><w:rPr><w:lang w:val="en-US"/></w:rPr></w:pPr></w:p>'
$p3.InsertXML($p3xml)

$p3again = $d.Paragraphs.Last.Range
$p3again.InsertParagraphAfter()
$p4 = $d.Paragraphs.Last.Range
$p4.InsertXML($p3xml)
